# The edit rotates the text content of 10 "value" paragraphs: each
# paragraph's text is replaced by the text that used to belong to the
# next paragraph in this chain (wrapping around at the end). Paragraph
# styles/formatting stay where they are; only the w:t content moves.
#
# Old text (in document order) -> New text
#  1. Complementar a formação...(PT objetivo)      -> A definir de acordo com o tópico programado
#  2. Complement the training...(EN objetivo)      -> To be defined according to the scheduled topic
#  3. 11079086 - Herlandí de Souza Andrade          -> Complementar a formação...(PT objetivo)
#  4. A definir de acordo com o tópico programado   -> O conteúdo desta disciplina optativa...
#  5. To be defined according to the scheduled topic-> Complement the training...(EN objetivo)
#  6. O conteúdo desta disciplina optativa...       -> Esta disciplina deverá conter...(Método)
#  7. Esta disciplina deverá conter...(Método)      -> Média ponderada das avaliações (M).
#  8. Média ponderada das avaliações (M).           -> A recuperação será composta...
#  9. A recuperação será composta...                -> Livros, artigos ou texto fornecido...
# 10. Livros, artigos ou texto fornecido...         -> 11079086 - Herlandí de Souza Andrade

$d = $word.ActiveDocument

# Use unique placeholder tokens first so the rotation doesn't clobber
# itself when one new value equals another old value's search target.
$pairs = @(
    @{ old = "Complementar a formação dos estudantes abordando, com maior profundidade, tópicos atuais e relevantes e atualizar com temas no estado da arte."; token = "@@TOKEN1@@" },
    @{ old = "Complement the training of students by addressing, in greater depth, current and relevant topics and updating them with themes in the state of the art"; token = "@@TOKEN2@@" },
    @{ old = "11079086 - Herlandí de Souza Andrade"; token = "@@TOKEN3@@" },
    @{ old = "A definir de acordo com o tópico programado"; token = "@@TOKEN4@@" },
    @{ old = "To be defined according to the scheduled topic"; token = "@@TOKEN5@@" },
    @{ old = "O conteúdo desta disciplina optativa será de acordo com o tópico a ser programado, devendo abordar assuntos complementares ao conteúdo regular do curso de graduação."; token = "@@TOKEN6@@" },
    @{ old = "Esta disciplina deverá conter no mínimo duas avaliações denominadas A1 e A2. As avalições poderão ser: escritas, práticas, seminários, trabalhos de campo, projetos, ou outra forma de avaliação definida pelo professor."; token = "@@TOKEN7@@" },
    @{ old = "Média ponderada das avaliações (M)."; token = "@@TOKEN8@@" },
    @{ old = "A recuperação será composta por uma única prova (RC) englobando toda a matéria ministrada ao longo do semestre. A média final, para os alunos em recuperação, será calculada com base na relação: MF=(M+RC)/2"; token = "@@TOKEN9@@" },
    @{ old = "Livros, artigos ou texto fornecido pelo docente responsável extraídos de livros ou revistas especializadas na área de Engenharia de Produção."; token = "@@TOKEN10@@" }
)

foreach ($pair in $pairs) {
    $d.Content.Find.Execute($pair.old, $true, $false, $false, $false, $false, $true, 1, $false, $pair.token, 2) | Out-Null
}

# token -> new text (the rotation)
$finals = @(
    @{ token = "@@TOKEN1@@"; new = "A definir de acordo com o tópico programado" },
    @{ token = "@@TOKEN2@@"; new = "To be defined according to the scheduled topic" },
    @{ token = "@@TOKEN3@@"; new = "Complementar a formação dos estudantes abordando, com maior profundidade, tópicos atuais e relevantes e atualizar com temas no estado da arte." },
    @{ token = "@@TOKEN4@@"; new = "O conteúdo desta disciplina optativa será de acordo com o tópico a ser programado, devendo abordar assuntos complementares ao conteúdo regular do curso de graduação." },
    @{ token = "@@TOKEN5@@"; new = "Complement the training of students by addressing, in greater depth, current and relevant topics and updating them with themes in the state of the art" },
    @{ token = "@@TOKEN6@@"; new = "Esta disciplina deverá conter no mínimo duas avaliações denominadas A1 e A2. As avalições poderão ser: escritas, práticas, seminários, trabalhos de campo, projetos, ou outra forma de avaliação definida pelo professor." },
    @{ token = "@@TOKEN7@@"; new = "Média ponderada das avaliações (M)." },
    @{ token = "@@TOKEN8@@"; new = "A recuperação será composta por uma única prova (RC) englobando toda a matéria ministrada ao longo do semestre. A média final, para os alunos em recuperação, será calculada com base na relação: MF=(M+RC)/2" },
    @{ token = "@@TOKEN9@@"; new = "Livros, artigos ou texto fornecido pelo docente responsável extraídos de livros ou revistas especializadas na área de Engenharia de Produção." },
    @{ token = "@@TOKEN10@@"; new = "11079086 - Herlandí de Souza Andrade" }
)

foreach ($pair in $finals) {
    $d.Content.Find.Execute($pair.token, $true, $false, $false, $false, $false, $true, 1, $false, $pair.new, 2) | Out-Null
}
